$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes existing rows 6..35 down to 7..36),
# mirroring the weekly price-update commit that added a new Madrigal entry.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 45063
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112013
$ws.Cells.Item(6, 7).Value = "Alcachofa"
$ws.Cells.Item(6, 8).Value = "Madrigal"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 21000
$ws.Cells.Item(6, 13).Value = 20500
$ws.Cells.Item(6, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(6, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(6, 16).Value = 512
$ws.Cells.Item(6, 17).Value = 40
$ws.Cells.Item(6, 18).Value = "Hortaliza"
